# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Islas Malvinas" / "Groenlandia" rows (A210 / A211) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 15:56"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 4102105
$ws.Range("C4").Value = 1230
$ws.Range("D4").Value = 1943504
$ws.Range("E4").Value = 2012402
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 146199

# --- Row 6: India ---
$ws.Range("B6").Value = 1257828
$ws.Range("C6").Value = 18144
$ws.Range("D6").Value = 796206
$ws.Range("E6").Value = 431564
$ws.Range("G6").Value = 168
$ws.Range("H6").Value = 30058

# --- Row 23: Argentina ---
$ws.Range("D23").Value = 62815
$ws.Range("E23").Value = 76468
$ws.Range("G23").Value = 29
$ws.Range("H23").Value = 2617

# --- Row 45: Portugal ---
$ws.Range("B45").Value = 49379
$ws.Range("C45").Value = 229
$ws.Range("D45").Value = 34369
$ws.Range("E45").Value = 13305
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 1705

# --- Row 58: Azerbaiyan ---
$ws.Range("B58").Value = 28980
$ws.Range("C58").Value = 347
$ws.Range("D58").Value = 20974
$ws.Range("E58").Value = 7615
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 391

# --- Row 62: Serbia ---
$ws.Range("B62").Value = 22443
$ws.Range("C62").Value = 412
$ws.Range("E62").Value = 7888
$ws.Range("G62").Value = 9
$ws.Range("H62").Value = 508

# --- Row 64: Austria ---
$ws.Range("B64").Value = 20099
$ws.Range("C64").Value = 170
$ws.Range("D64").Value = 17943
$ws.Range("E64").Value = 1445

# --- Row 80: Estado de Palestina ---
$ws.Range("E80").Value = 6957
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 67

# --- Row 85: Noruega ---
$ws.Range("B85").Value = 9071
$ws.Range("C85").Value = 12
$ws.Range("E85").Value = 142

# --- Row 91: Tayikistan ---
$ws.Range("B91").Value = 7060
$ws.Range("C91").Value = 45
$ws.Range("D91").Value = 5793
$ws.Range("E91").Value = 1209

# --- Row 117: Cuba ---
$ws.Range("B117").Value = 2466
$ws.Range("C117").Value = 4
$ws.Range("D117").Value = 2339
$ws.Range("E117").Value = 40

# --- Row 131: Benin ---
$ws.Range("B131").Value = 1694
$ws.Range("C131").Value = 4
$ws.Range("D131").Value = 918
$ws.Range("E131").Value = 742

# --- Row 145: Uganda ---
$ws.Range("B145").Value = 1079
$ws.Range("C145").Value = 4
$ws.Range("D145").Value = 971
$ws.Range("E145").Value = 108

# --- Row 146: Burkina Faso ---
$ws.Range("B146").Value = 1070
$ws.Range("C146").Value = 4
$ws.Range("D146").Value = 919
$ws.Range("E146").Value = 98
